# Applies the update described by the commit:
# "Atualizado por script em 27-10-2023 02:45"
#
# 1) Swap the betting-data (columns F:V) of rows 39 <-> 40
# 2) Swap the betting-data (columns F:V) of rows 43 <-> 44
# 3) Append 5 new match rows (176-180) with freshly scraped odds

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap rows 39 and 40 (only columns F:V -- home/away teams,
#    scores, odds, timestamps and url; Indice/pais/torneio/
#    temporada/data_partida in A:E stay untouched per row)
# ---------------------------------------------------------------
$row39 = $ws.Range("F39:V39")
$row40 = $ws.Range("F40:V40")

$tmp39 = $row39.Value2
$tmp40 = $row40.Value2

$row39.Value2 = $tmp40
$row40.Value2 = $tmp39

# ---------------------------------------------------------------
# 2) Swap rows 43 and 44 (same rule as above)
# ---------------------------------------------------------------
$row43 = $ws.Range("F43:V43")
$row44 = $ws.Range("F44:V44")

$tmp43 = $row43.Value2
$tmp44 = $row44.Value2

$row43.Value2 = $tmp44
$row44.Value2 = $tmp43

# ---------------------------------------------------------------
# 3) Append new rows 176-180.
#    Each new row is seeded by copying the formatting (and, for
#    the shared columns B/C/D = ecuador/serie-b/2023, the values
#    too) from row 175, then the per-match fields are overwritten.
# ---------------------------------------------------------------

# --- New row 176 (Indice=175) ---
$src = $ws.Range("A175:V175")
$dst = $ws.Range("A176:V176")
$src.Copy($dst)
$ws.Cells.Item(176,1).Value2 = 175
$ws.Cells.Item(176,5).Value2 = 45226.08333333334
$ws.Cells.Item(176,6).Value2 = "Ind. Juniors"
$ws.Cells.Item(176,7).Value2 = 2
$ws.Cells.Item(176,8).Value2 = "Cuniburo"
$ws.Cells.Item(176,9).Value2 = 4
$ws.Cells.Item(176,10).Value2 = 1.77
$ws.Cells.Item(176,11).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(176,12).Value2 = 3.36
$ws.Cells.Item(176,13).Value2 = "27/10/2023 01:58"
$ws.Cells.Item(176,14).Value2 = 3.47
$ws.Cells.Item(176,15).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(176,16).Value2 = 3.82
$ws.Cells.Item(176,17).Value2 = "27/10/2023 01:58"
$ws.Cells.Item(176,18).Value2 = 3.95
$ws.Cells.Item(176,19).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(176,20).Value2 = 1.95
$ws.Cells.Item(176,21).Value2 = "27/10/2023 01:58"
$ws.Cells.Item(176,22).Value2 = "https://www.betexplorer.com/football/ecuador/serie-b/independiente-juniors-cuniburo/MithnJF8/"

# --- New row 177 (Indice=176) ---
$src = $ws.Range("A175:V175")
$dst = $ws.Range("A177:V177")
$src.Copy($dst)
$ws.Cells.Item(177,1).Value2 = 176
$ws.Cells.Item(177,5).Value2 = 45226.08333333334
$ws.Cells.Item(177,6).Value2 = "America de Quito"
$ws.Cells.Item(177,7).Value2 = 3
$ws.Cells.Item(177,8).Value2 = "Chacaritas"
$ws.Cells.Item(177,9).Value2 = 1
$ws.Cells.Item(177,10).Value2 = 2.2
$ws.Cells.Item(177,11).Value2 = "25/10/2023 13:15"
$ws.Cells.Item(177,12).Value2 = 1.13
$ws.Cells.Item(177,13).Value2 = "27/10/2023 01:59"
$ws.Cells.Item(177,14).Value2 = 3.07
$ws.Cells.Item(177,15).Value2 = "25/10/2023 13:15"
$ws.Cells.Item(177,16).Value2 = 8.449999999999999
$ws.Cells.Item(177,17).Value2 = "27/10/2023 01:59"
$ws.Cells.Item(177,18).Value2 = 3.08
$ws.Cells.Item(177,19).Value2 = "25/10/2023 13:15"
$ws.Cells.Item(177,20).Value2 = 14.73
$ws.Cells.Item(177,21).Value2 = "27/10/2023 01:59"
$ws.Cells.Item(177,22).Value2 = "https://www.betexplorer.com/football/ecuador/serie-b/america-de-quito-chacaritas/GEudowVE/"

# --- New row 178 (Indice=177) ---
$src = $ws.Range("A175:V175")
$dst = $ws.Range("A178:V178")
$src.Copy($dst)
$ws.Cells.Item(178,1).Value2 = 177
$ws.Cells.Item(178,5).Value2 = 45226.08333333334
$ws.Cells.Item(178,6).Value2 = "Buhos ULVR"
$ws.Cells.Item(178,7).Value2 = 0
$ws.Cells.Item(178,8).Value2 = "Vargas Torres"
$ws.Cells.Item(178,9).Value2 = 1
$ws.Cells.Item(178,10).Value2 = 2.15
$ws.Cells.Item(178,11).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(178,12).Value2 = 5.66
$ws.Cells.Item(178,13).Value2 = "27/10/2023 01:53"
$ws.Cells.Item(178,14).Value2 = 3.08
$ws.Cells.Item(178,15).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(178,16).Value2 = 3.99
$ws.Cells.Item(178,17).Value2 = "27/10/2023 01:53"
$ws.Cells.Item(178,18).Value2 = 3.18
$ws.Cells.Item(178,19).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(178,20).Value2 = 1.55
$ws.Cells.Item(178,21).Value2 = "27/10/2023 01:53"
$ws.Cells.Item(178,22).Value2 = "https://www.betexplorer.com/football/ecuador/serie-b/buhos-ulvr-vargas-torres/zHecqHaR/"

# --- New row 179 (Indice=178) ---
$src = $ws.Range("A175:V175")
$dst = $ws.Range("A179:V179")
$src.Copy($dst)
$ws.Cells.Item(179,1).Value2 = 178
$ws.Cells.Item(179,5).Value2 = 45226.08333333334
$ws.Cells.Item(179,6).Value2 = "Imbabura"
$ws.Cells.Item(179,7).Value2 = 4
$ws.Cells.Item(179,8).Value2 = "Macara"
$ws.Cells.Item(179,9).Value2 = 0
$ws.Cells.Item(179,10).Value2 = 2.24
$ws.Cells.Item(179,11).Value2 = "25/10/2023 13:15"
$ws.Cells.Item(179,12).Value2 = 1.96
$ws.Cells.Item(179,13).Value2 = "27/10/2023 01:59"
$ws.Cells.Item(179,14).Value2 = 3
$ws.Cells.Item(179,15).Value2 = "25/10/2023 13:15"
$ws.Cells.Item(179,16).Value2 = 2.7
$ws.Cells.Item(179,17).Value2 = "27/10/2023 01:59"
$ws.Cells.Item(179,18).Value2 = 3.09
$ws.Cells.Item(179,19).Value2 = "25/10/2023 13:15"
$ws.Cells.Item(179,20).Value2 = 5.25
$ws.Cells.Item(179,21).Value2 = "27/10/2023 01:59"
$ws.Cells.Item(179,22).Value2 = "https://www.betexplorer.com/football/ecuador/serie-b/imbabura-macara/f5v0pcpL/"

# --- New row 180 (Indice=179) ---
$src = $ws.Range("A175:V175")
$dst = $ws.Range("A180:V180")
$src.Copy($dst)
$ws.Cells.Item(180,1).Value2 = 179
$ws.Cells.Item(180,5).Value2 = 45226.08333333334
$ws.Cells.Item(180,6).Value2 = "Manta"
$ws.Cells.Item(180,7).Value2 = 0
$ws.Cells.Item(180,8).Value2 = "Nueve de Octubre"
$ws.Cells.Item(180,9).Value2 = 1
$ws.Cells.Item(180,10).Value2 = 2.02
$ws.Cells.Item(180,11).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(180,12).Value2 = 3.96
$ws.Cells.Item(180,13).Value2 = "27/10/2023 01:58"
$ws.Cells.Item(180,14).Value2 = 3.14
$ws.Cells.Item(180,15).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(180,16).Value2 = 4
$ws.Cells.Item(180,17).Value2 = "27/10/2023 01:58"
$ws.Cells.Item(180,18).Value2 = 3.42
$ws.Cells.Item(180,19).Value2 = "25/10/2023 13:14"
$ws.Cells.Item(180,20).Value2 = 1.76
$ws.Cells.Item(180,21).Value2 = "27/10/2023 01:58"
$ws.Cells.Item(180,22).Value2 = "https://www.betexplorer.com/football/ecuador/serie-b/manta-nueve-de-octubre/xQsPvy8r/"

Write-Host "Update complete: swapped rows 39/40, 43/44 and appended rows 176-180."
